$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (shifts existing rows 7-13 down to 8-14)
$ws.Rows.Item(7).Insert()

# Fill in the new census (Censos) row with the latest census data point
$ws.Range("A7").Value = "Censos"
$ws.Range("B7").Value = 2023
$ws.Range("C7").Value = 3444.2629999999999

# Add the new projection (Proy) row at the end of the table
$ws.Range("A15").Value = "Proy (NNUU 2019) "
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = 3473.7269999999999

# Update the filter database named range to include the newly added rows
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Hoja1!`$A`$1:`$C`$14"

# Re-apply the sort over the expanded data range so the sort state reflects the new rows
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A14"))
$ws.Sort.SortFields.Add($ws.Range("B2:B14"))
$ws.Sort.SetRange($ws.Range("A2:C14"))
$ws.Sort.Apply()

# Auto-fit column A so the longer "Proy (NNUU 2019)" label is fully visible
$ws.Columns.Item(1).AutoFit()

# Match the active selection left after the edit
$ws.Range("C8").Select() | Out-Null

$wb.Save() | Out-Null
